$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 314 (shifts rows 314-338 down to 315-339)
$ws.Rows.Item(314).Insert()

# Populate the newly inserted row 314 with the new record
$ws.Range("A314").Value = 5
$ws.Range("B314").Value = "Macroferia Regional de Talca"
$ws.Range("C314").Value = "Maule"
$ws.Range("D314").Value2 = 44783
$ws.Range("E314").Value = 7
$ws.Range("F314").Value = 100112003
$ws.Range("G314").Value = "Ajo"
$ws.Range("H314").Value = "Chino"
$ws.Range("I314").Value = "Primera"
$ws.Range("J314").Value = 300
$ws.Range("K314").Value = 25000
$ws.Range("L314").Value = 25000
$ws.Range("M314").Value = 25000
$ws.Range("N314").Value = "`$/caja 10 kilos"
$ws.Range("O314").Value = "China"
$ws.Range("P314").Value = 2500
$ws.Range("Q314").Value = 10
$ws.Range("R314").Value = "Hortaliza"
